$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44270
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 85
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("R2").Value = 'Provincia del Elquí'
$ws.Range("S2").Value = 857

$ws.Range("D3").Value = 44245
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = '$/caja 15 kilos granel'
$ws.Range("T3").Value = 15

$ws.Range("D4").Value = 44278
$ws.Range("L4").Value = 'Primera'
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 13000
$ws.Range("Q4").Value = '$/caja 14 kilos empedrada'
$ws.Range("R4").Value = 'Provincia del Elquí'
$ws.Range("S4").Value = 929

$ws.Range("D5").Value = 45006
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("Q5").Value = '$/caja 14 kilos empedrada'
$ws.Range("S5").Value = 1143

$ws.Range("D6").Value = 44315
$ws.Range("M6").Value = 65

$ws.Range("D7").Value = 45015
$ws.Range("M7").Value = 56
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("Q7").Value = '$/caja 14 kilos empedrada'
$ws.Range("S7").Value = 1071
$ws.Range("T7").Value = 14

$ws.Range("D8").Value = 44592
$ws.Range("M8").Value = 54
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("Q8").Value = '$/caja 15 kilos empedrada'
$ws.Range("S8").Value = 1333
$ws.Range("T8").Value = 15

$ws.Range("D9").Value = 44614
$ws.Range("M9").Value = 54
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 14000
$ws.Range("Q9").Value = '$/caja 14 kilos granel'
$ws.Range("S9").Value = 1000

$ws.Range("D10").Value = 44314
$ws.Range("M10").Value = 56
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 14000
$ws.Range("S10").Value = 1000

$ws.Range("D11").Value = 45014
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("Q11").Value = '$/caja 14 kilos empedrada'
$ws.Range("S11").Value = 1071

$ws.Range("D12").Value = 44239
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("Q12").Value = '$/caja 15 kilos granel'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 1000
$ws.Range("T12").Value = 15

$ws.Range("D13").Value = 44323
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 14000
$ws.Range("Q13").Value = '$/caja 14 kilos granel'
$ws.Range("T13").Value = 14

$ws.Range("D14").Value = 44630
$ws.Range("M14").Value = 75
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("S14").Value = 1071

$ws.Range("D15").Value = 45040
$ws.Range("L15").Value = 'Especial'
$ws.Range("M15").Value = 65
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 17000
$ws.Range("P15").Value = 17000
$ws.Range("Q15").Value = '$/caja 14 kilos granel'
$ws.Range("S15").Value = 1214

$ws.Range("D16").Value = 45040
$ws.Range("M16").Value = 60

$ws.Range("D17").Value = 45054
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 54
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("Q17").Value = '$/caja 14 kilos empedrada'
$ws.Range("S17").Value = 1143

$ws.Range("D18").Value = 45054
$ws.Range("M18").Value = 50

$ws.Range("D19").Value = 44320
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 45

$ws.Range("L20").Value = 'Especial'
$ws.Range("M20").Value = 56
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 14000
$ws.Range("S20").Value = 1000

$ws.Range("D21").Value = 45050
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 12000
$ws.Range("Q21").Value = '$/caja 14 kilos granel'
$ws.Range("S21").Value = 857

$ws.Range("D22").Value = 44312
$ws.Range("M22").Value = 68
$ws.Range("N22").Value = 14000
$ws.Range("O22").Value = 14000
$ws.Range("P22").Value = 14000
$ws.Range("Q22").Value = '$/caja 14 kilos granel'
$ws.Range("R22").Value = 'Provincia de Limarí'
$ws.Range("S22").Value = 1000

$ws.Range("D23").Value = 44242
$ws.Range("M23").Value = 45
$ws.Range("N23").Value = 12000
$ws.Range("O23").Value = 12000
$ws.Range("P23").Value = 12000
$ws.Range("S23").Value = 800

$ws.Range("D24").Value = 44260
$ws.Range("M24").Value = 56
$ws.Range("N24").Value = 13000
$ws.Range("O24").Value = 13000
$ws.Range("P24").Value = 13000
$ws.Range("Q24").Value = '$/caja 14 kilos empedrada'
$ws.Range("R24").Value = 'Provincia del Elquí'
$ws.Range("S24").Value = 929

$ws.Range("D25").Value = 44259
$ws.Range("M25").Value = 80
$ws.Range("N25").Value = 12000
$ws.Range("O25").Value = 12000
$ws.Range("P25").Value = 12000
$ws.Range("Q25").Value = '$/caja 15 kilos empedrada'
$ws.Range("S25").Value = 800
$ws.Range("T25").Value = 15

$ws.Range("D26").Value = 45044
$ws.Range("M26").Value = 30
$ws.Range("N26").Value = 16000
$ws.Range("O26").Value = 16000
$ws.Range("P26").Value = 16000
$ws.Range("S26").Value = 1143

$ws.Range("D27").Value = 45044
$ws.Range("M27").Value = 30

$ws.Range("D29").Value = 44322
$ws.Range("N29").Value = 14000
$ws.Range("O29").Value = 14000
$ws.Range("P29").Value = 14000
$ws.Range("R29").Value = 'Provincia de Limarí'
$ws.Range("S29").Value = 1000

$ws.Range("D30").Value = 44313
$ws.Range("M30").Value = 36

$ws.Range("D31").Value = 44271
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 12000
$ws.Range("O31").Value = 12000
$ws.Range("P31").Value = 12000
$ws.Range("R31").Value = 'Provincia del Elquí'
$ws.Range("S31").Value = 857

$ws.Range("D32").Value = 44252
$ws.Range("M32").Value = 60
$ws.Range("N32").Value = 14000
$ws.Range("O32").Value = 14000
$ws.Range("P32").Value = 14000
$ws.Range("Q32").Value = '$/caja 14 kilos empedrada'
$ws.Range("R32").Value = 'Provincia de Limarí'
$ws.Range("S32").Value = 1000

$ws.Range("D33").Value = 44588
$ws.Range("M33").Value = 85
$ws.Range("N33").Value = 19000
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 19529
$ws.Range("Q33").Value = '$/caja 14 kilos granel'
$ws.Range("S33").Value = 1395

$ws.Range("D36").Value = 44627
$ws.Range("M36").Value = 56
$ws.Range("N36").Value = 17000
$ws.Range("O36").Value = 17000
$ws.Range("P36").Value = 17000
$ws.Range("Q36").Value = '$/caja 14 kilos empedrada'
$ws.Range("S36").Value = 1214
$ws.Range("T36").Value = 14

$ws.Range("D37").Value = 45043
$ws.Range("M37").Value = 45
$ws.Range("N37").Value = 17000
$ws.Range("O37").Value = 17000
$ws.Range("P37").Value = 17000
$ws.Range("Q37").Value = '$/caja 14 kilos granel'
$ws.Range("S37").Value = 1214

$ws.Range("D38").Value = 45043
$ws.Range("M38").Value = 67
$ws.Range("Q38").Value = '$/caja 14 kilos granel'

$ws.Range("D39").Value = 44316
$ws.Range("M39").Value = 48
$ws.Range("N39").Value = 14000
$ws.Range("O39").Value = 14000
$ws.Range("P39").Value = 14000
$ws.Range("Q39").Value = '$/caja 14 kilos granel'
$ws.Range("S39").Value = 1000

$ws.Range("D40").Value = 44238
$ws.Range("M40").Value = 60
$ws.Range("N40").Value = 15000
$ws.Range("O40").Value = 15000
$ws.Range("P40").Value = 15000
$ws.Range("Q40").Value = '$/caja 15 kilos granel'
$ws.Range("S40").Value = 1000

$ws.Range("D41").Value = 45001
$ws.Range("M41").Value = 50
$ws.Range("N41").Value = 16000
$ws.Range("O41").Value = 16000
$ws.Range("P41").Value = 16000
$ws.Range("Q41").Value = '$/caja 14 kilos empedrada'
$ws.Range("S41").Value = 1143
$ws.Range("T41").Value = 14

